# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (per scheduled price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4055.889
$ws.Range("I74").Value = 4460
$ws.Range("J74").Value = 3550.75
$ws.Range("K74").Value = 4460
$ws.Range("L74").Value = 3550.75
$ws.Range("M74").Value = -3524
$ws.Range("N74").Value = -5422.75
# Row 77
$ws.Range("H77").Value = 4055.889
$ws.Range("I77").Value = 4460
$ws.Range("J77").Value = 3550.75
$ws.Range("K77").Value = 22300
$ws.Range("L77").Value = 17753.75
$ws.Range("M77").Value = -17620
$ws.Range("N77").Value = -27113.75
# Row 99
$ws.Range("H99").Value = 4040.75
$ws.Range("I99").Value = 4040.75
$ws.Range("K99").Value = 12122.25
$ws.Range("M99").Value = -10624.25
# Row 109
$ws.Range("H109").Value = 38263
$ws.Range("J109").Value = 38263
$ws.Range("L109").Value = 38263
$ws.Range("N109").Value = -41037
# Row 117
$ws.Range("H117").Value = 46178
$ws.Range("J117").Value = 46178
$ws.Range("L117").Value = 46178
$ws.Range("N117").Value = -55356
# Row 124
$ws.Range("H124").Value = 41992
$ws.Range("J124").Value = 41992
$ws.Range("L124").Value = 41992
$ws.Range("N124").Value = -51812
# Row 128
$ws.Range("H128").Value = 45722.8
$ws.Range("J128").Value = 45722.8
$ws.Range("L128").Value = 45722.8
$ws.Range("N128").Value = -55682.8
# Row 130
$ws.Range("H130").Value = 46169
$ws.Range("J130").Value = 46169
$ws.Range("L130").Value = 46169
$ws.Range("N130").Value = -56209

$ws = $wb.Worksheets.Item("ARM")
# Row 113
$ws.Range("H113").Value = 41266.57
$ws.Range("J113").Value = 41266.57
$ws.Range("L113").Value = 41266.57
$ws.Range("N113").Value = -49944.57
# Row 114
$ws.Range("H114").Value = 44961.332
$ws.Range("J114").Value = 44961.332
$ws.Range("L114").Value = 44961.332
$ws.Range("N114").Value = -53639.332
# Row 118
$ws.Range("H118").Value = 43281.832
$ws.Range("J118").Value = 43281.832
$ws.Range("L118").Value = 43281.832
$ws.Range("N118").Value = -46595.832
# Row 123
$ws.Range("H123").Value = 49992
$ws.Range("J123").Value = 49992
$ws.Range("L123").Value = 49992
$ws.Range("N123").Value = -59792
# Row 125
$ws.Range("H125").Value = 48740.668
$ws.Range("J125").Value = 48740.668
$ws.Range("L125").Value = 48740.668
$ws.Range("N125").Value = -58580.668
# Row 130
$ws.Range("H130").Value = 38483.332
$ws.Range("J130").Value = 38483.332
$ws.Range("L130").Value = 38483.332
$ws.Range("N130").Value = -48523.332
# Row 131
$ws.Range("H131").Value = 51711
$ws.Range("J131").Value = 51711
$ws.Range("L131").Value = 51711
$ws.Range("N131").Value = -61791

$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Range("H108").Value = 46070.8
$ws.Range("J108").Value = 46070.8
$ws.Range("L108").Value = 46070.8
$ws.Range("N108").Value = -53750.8
# Row 111
$ws.Range("H111").Value = 37900.668
$ws.Range("J111").Value = 37900.668
$ws.Range("L111").Value = 37900.668
$ws.Range("N111").Value = -46080.668
# Row 117
$ws.Range("H117").Value = 49248
$ws.Range("J117").Value = 49248
$ws.Range("L117").Value = 49248
$ws.Range("N117").Value = -58426
# Row 124
$ws.Range("H124").Value = 49881.332
$ws.Range("J124").Value = 49881.332
$ws.Range("L124").Value = 49881.332
$ws.Range("N124").Value = -59701.332
# Row 125
$ws.Range("H125").Value = 50472
$ws.Range("J125").Value = 50472
$ws.Range("L125").Value = 50472
$ws.Range("N125").Value = -60312
# Row 126
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652
# Row 130
$ws.Range("H130").Value = 49514
$ws.Range("J130").Value = 49514
$ws.Range("L130").Value = 49514
$ws.Range("N130").Value = -59554

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 40664.5
$ws.Range("J20").Value = 40664.5
$ws.Range("L20").Value = 40664.5
$ws.Range("N20").Value = -41136.5
# Row 30
$ws.Range("H30").Value = 40664.5
$ws.Range("J30").Value = 40664.5
$ws.Range("L30").Value = 40664.5
$ws.Range("N30").Value = -40846.5
# Row 116
$ws.Range("H116").Value = 47822.332
$ws.Range("J116").Value = 47822.332
$ws.Range("L116").Value = 47822.332
$ws.Range("N116").Value = -57000.332
# Row 119
$ws.Range("H119").Value = 48261
$ws.Range("J119").Value = 48261
$ws.Range("L119").Value = 48261
$ws.Range("N119").Value = -57937
# Row 128
$ws.Range("H128").Value = 40664.5
$ws.Range("J128").Value = 40664.5
$ws.Range("L128").Value = 40664.5
$ws.Range("N128").Value = -50624.5

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1508.9474
$ws.Range("I132").Value = 1076.2222
$ws.Range("J132").Value = 1898.4
$ws.Range("K132").Value = 9685.9998
$ws.Range("L132").Value = 17085.6
$ws.Range("M132").Value = -7155.9998
$ws.Range("N132").Value = -22145.6

$ws = $wb.Worksheets.Item("GSM")
# Row 110
$ws.Range("H110").Value = 48702
$ws.Range("J110").Value = 48702
$ws.Range("L110").Value = 48702
$ws.Range("N110").Value = -56882
# Row 116
$ws.Range("H116").Value = 49434
$ws.Range("J116").Value = 49434
$ws.Range("L116").Value = 49434
$ws.Range("N116").Value = -58612
# Row 122
$ws.Range("H122").Value = 2224.2856
$ws.Range("I122").Value = 2556.6667
$ws.Range("J122").Value = 1975
$ws.Range("K122").Value = 7670.000100000001
$ws.Range("L122").Value = 5925
$ws.Range("M122").Value = -5220.000100000001
$ws.Range("N122").Value = -10825
# Row 130
$ws.Range("H130").Value = 44728
$ws.Range("J130").Value = 44728
$ws.Range("L130").Value = 44728
$ws.Range("N130").Value = -54768

$ws = $wb.Worksheets.Item("LTW")
# Row 114
$ws.Range("H114").Value = 38344
$ws.Range("J114").Value = 38344
$ws.Range("L114").Value = 38344
$ws.Range("N114").Value = -47022
# Row 116
$ws.Range("H116").Value = 50670.2
$ws.Range("J116").Value = 50670.2
$ws.Range("L116").Value = 50670.2
$ws.Range("N116").Value = -59848.2
# Row 124
$ws.Range("H124").Value = 43140.332
$ws.Range("J124").Value = 43140.332
$ws.Range("L124").Value = 43140.332
$ws.Range("N124").Value = -52960.332
# Row 125
$ws.Range("H125").Value = 49715
$ws.Range("J125").Value = 49715
$ws.Range("L125").Value = 49715
$ws.Range("N125").Value = -59555
# Row 127
$ws.Range("H127").Value = 50496
$ws.Range("J127").Value = 50496
$ws.Range("L127").Value = 50496
$ws.Range("N127").Value = -60416
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 130
$ws.Range("H130").Value = 41723.11
$ws.Range("J130").Value = 41723.11
$ws.Range("L130").Value = 41723.11
$ws.Range("N130").Value = -51763.11

$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 48618
$ws.Range("J108").Value = 48618
$ws.Range("L108").Value = 48618
$ws.Range("N108").Value = -56298
# Row 110
$ws.Range("H110").Value = 46994.668
$ws.Range("J110").Value = 46994.668
$ws.Range("L110").Value = 46994.668
$ws.Range("N110").Value = -55174.668
# Row 116
$ws.Range("H116").Value = 47657.332
$ws.Range("J116").Value = 47657.332
$ws.Range("L116").Value = 47657.332
$ws.Range("N116").Value = -56835.332
# Row 120
$ws.Range("H120").Value = 42460.8
$ws.Range("J120").Value = 42460.8
$ws.Range("L120").Value = 42460.8
$ws.Range("N120").Value = -52136.8
# Row 121
$ws.Range("H121").Value = 33975.2
$ws.Range("J121").Value = 33975.2
$ws.Range("L121").Value = 33975.2
$ws.Range("N121").Value = -37469.2
# Row 128
$ws.Range("H128").Value = 50711
$ws.Range("J128").Value = 50711
$ws.Range("L128").Value = 50711
$ws.Range("N128").Value = -60671
# Row 131
$ws.Range("H131").Value = 49232.25
$ws.Range("J131").Value = 49232.25
$ws.Range("L131").Value = 49232.25
$ws.Range("N131").Value = -59312.25

Write-Host "Applied Leve profit updates."
